$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("K2").Value = 4392
$ws.Range("K3").Value = 4497
$ws.Range("B4").Value = 1703
$ws.Range("K4").Value = 901
$ws.Range("K5").Value = 325
$ws.Range("K6").Value = 5058
$ws.Range("B7").Value = 23336
$ws.Range("K7").Value = 15173

$ws = $wb.Worksheets.Item('Grant Park')
$ws.Range("K3").Value = 5
$ws.Range("K5").Value = 6
$ws.Range("K6").Value = 14

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("K3").Value = 48
$ws.Range("K7").Value = 202

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("K2").Value = 283
$ws.Range("J4").Value = 97
$ws.Range("K6").Value = 342
$ws.Range("J7").Value = 1852
$ws.Range("K7").Value = 1013

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("K2").Value = 112
$ws.Range("K3").Value = 117
$ws.Range("K6").Value = 73
$ws.Range("K7").Value = 325

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("K2").Value = 175
$ws.Range("K3").Value = 236
$ws.Range("K6").Value = 184
$ws.Range("K7").Value = 637

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("K2").Value = 89
$ws.Range("K6").Value = 65
$ws.Range("K7").Value = 269

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("K2").Value = 143
$ws.Range("K6").Value = 154
$ws.Range("K7").Value = 514

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("K3").Value = 106
$ws.Range("K7").Value = 252

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("K2").Value = 132
$ws.Range("K4").Value = 58
$ws.Range("K7").Value = 445
$ws.Range("J8").Value = 1852
$ws.Range("K8").Value = 1013
$ws.Range("K11").Value = 300
$ws.Range("K15").Value = 153
$ws.Range("K16").Value = 51
$ws.Range("K17").Value = 28
$ws.Range("K19").Value = 460
$ws.Range("K20").Value = 345
$ws.Range("K23").Value = 156
$ws.Range("J27").Value = 180
$ws.Range("K29").Value = 802
$ws.Range("K33").Value = 637
$ws.Range("K35").Value = 22
$ws.Range("K37").Value = 514
$ws.Range("K38").Value = 14
$ws.Range("K39").Value = 21
$ws.Range("K41").Value = 121
$ws.Range("K42").Value = 554
$ws.Range("K43").Value = 137
$ws.Range("K48").Value = 198
$ws.Range("K50").Value = 82
$ws.Range("K51").Value = 194
$ws.Range("K52").Value = 404
$ws.Range("K53").Value = 202
$ws.Range("K54").Value = 284
$ws.Range("K55").Value = 170
$ws.Range("K57").Value = 52
$ws.Range("B63").Value = 407
$ws.Range("K63").Value = 49
$ws.Range("K67").Value = 585
$ws.Range("K71").Value = 49
$ws.Range("K73").Value = 134
$ws.Range("K76").Value = 208
$ws.Range("K78").Value = 178
$ws.Range("K80").Value = 51
$ws.Range("K83").Value = 325
$ws.Range("K85").Value = 679
$ws.Range("K88").Value = 177
$ws.Range("K89").Value = 218
$ws.Range("K93").Value = 55
$ws.Range("K94").Value = 188
$ws.Range("K95").Value = 269
$ws.Range("K96").Value = 168
$ws.Range("K99").Value = 252
$ws.Range("B101").Value = 23336
$ws.Range("K101").Value = 15173

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("K2").Value = 168
$ws.Range("K3").Value = 205
$ws.Range("K6").Value = 171
$ws.Range("K7").Value = 585

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("K3").Value = 79
$ws.Range("K7").Value = 284

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("K2").Value = 229
$ws.Range("K3").Value = 286
$ws.Range("K6").Value = 223
$ws.Range("K7").Value = 802

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("K4").Value = 27
$ws.Range("K6").Value = 97
$ws.Range("K7").Value = 198

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("K3").Value = 140
$ws.Range("K6").Value = 144
$ws.Range("K7").Value = 460

$ws = $wb.Worksheets.Item('River North')
$ws.Range("K6").Value = 116
$ws.Range("K7").Value = 208

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range("K2").Value = 43
$ws.Range("K7").Value = 121

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("K6").Value = 202
$ws.Range("K7").Value = 554

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("K3").Value = 39
$ws.Range("K7").Value = 178

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("K6").Value = 64
$ws.Range("K7").Value = 170

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("K2").Value = 44
$ws.Range("K3").Value = 58
$ws.Range("K6").Value = 39
$ws.Range("K7").Value = 156

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("K6").Value = 75
$ws.Range("K7").Value = 168

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("K2").Value = 120
$ws.Range("K3").Value = 107
$ws.Range("K7").Value = 345

$ws = $wb.Worksheets.Item('Burnside')
$ws.Range("K6").Value = 7
$ws.Range("K7").Value = 28

$ws = $wb.Worksheets.Item('West Lawn')
$ws.Range("K6").Value = 21
$ws.Range("K7").Value = 55

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("K5").Value = 19
$ws.Range("K6").Value = 114
$ws.Range("K7").Value = 445

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("K2").Value = 55
$ws.Range("K7").Value = 188

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("K2").Value = 53
$ws.Range("K7").Value = 153

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range("K3").Value = 12
$ws.Range("K7").Value = 82

$ws = $wb.Worksheets.Item('Greektown')
$ws.Range("K5").Value = 12
$ws.Range("K6").Value = 21

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("K3").Value = 75
$ws.Range("K7").Value = 300

$ws = $wb.Worksheets.Item('Gold Coast')
$ws.Range("K6").Value = 15
$ws.Range("K7").Value = 22

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range("K4").Value = 9
$ws.Range("K7").Value = 134

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("K2").Value = 42
$ws.Range("K7").Value = 132

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("K3").Value = 57
$ws.Range("K6").Value = 73
$ws.Range("K7").Value = 177

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("K2").Value = 57
$ws.Range("K3").Value = 66
$ws.Range("K7").Value = 218

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("J4").Value = 22
$ws.Range("J7").Value = 180

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("K2").Value = 49
$ws.Range("K6").Value = 28

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("K3").Value = 55
$ws.Range("K6").Value = 66
$ws.Range("K7").Value = 194

$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Range("K6").Value = 26
$ws.Range("K7").Value = 52

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("K6").Value = 58
$ws.Range("K7").Value = 137

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("K2").Value = 238
$ws.Range("K6").Value = 160
$ws.Range("K7").Value = 679

$ws = $wb.Worksheets.Item('Oakland')
$ws.Range("K3").Value = 19
$ws.Range("K7").Value = 49

$ws = $wb.Worksheets.Item('Rush & Division')
$ws.Range("K6").Value = 26
$ws.Range("K7").Value = 51

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("K3").Value = 108
$ws.Range("K4").Value = 24
$ws.Range("K7").Value = 404

$ws = $wb.Worksheets.Item('Archer Heights')
$ws.Range("K2").Value = 19
$ws.Range("K7").Value = 58

$ws = $wb.Worksheets.Item('Bucktown')
$ws.Range("K6").Value = 29
$ws.Range("K7").Value = 51
